$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.602.64"
$ws.Range("E2").Value = "  -0.04%  "

Set-TextValue $ws.Range("D3") "1.719.31"
$ws.Range("E3").Value = "  -1.14%  "

Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  +0.17%  "

Set-TextValue $ws.Range("D5") "240.89"
$ws.Range("E5").Value = "  -2.15%  "

Set-TextValue $ws.Range("D6") "1.000"
$ws.Range("E6").Value = "  +0.04%  "

Set-TextValue $ws.Range("D7") "0.4911"
$ws.Range("E7").Value = "  -0.38%  "

Set-TextValue $ws.Range("D8") "0.2597"
$ws.Range("E8").Value = "  -2.81%  "

Set-TextValue $ws.Range("D9") "0.06188"
$ws.Range("E9").Value = "  -1.37%  "

Set-TextValue $ws.Range("D10") "1.731.74"
$ws.Range("E10").Value = "  -1.01%  "

Set-TextValue $ws.Range("D11") "0.06977"
$ws.Range("E11").Value = "  -0.97%  "

Set-TextValue $ws.Range("D12") "15.65"
$ws.Range("E12").Value = "  -0.65%  "

Set-TextValue $ws.Range("D13") "0.6066"
$ws.Range("E13").Value = "  -1.18%  "

Set-TextValue $ws.Range("D14") "4.464"
$ws.Range("E14").Value = "  -2.52%  "

Set-TextValue $ws.Range("D15") "76.72"
$ws.Range("E15").Value = "  -1.64%  "

Set-TextValue $ws.Range("D16") "0.9975"
$ws.Range("E16").Value = "  -0.27%  "

Set-TextValue $ws.Range("D17") "26.459.81"
$ws.Range("E17").Value = "  -0.59%  "

Set-TextValue $ws.Range("D18") "1.003"
$ws.Range("E18").Value = "  +0.26%  "

Set-TextValue $ws.Range("D19") "0.000007122"
$ws.Range("E19").Value = "  -1.95%  "

Set-TextValue $ws.Range("D20") "11.34"
$ws.Range("E20").Value = "  -1.93%  "

Set-TextValue $ws.Range("D21") "1.953.28"
$ws.Range("E21").Value = "  -0.90%  "

Set-TextValue $ws.Range("D22") "4.413"
$ws.Range("E22").Value = "  -3.31%  "

Set-TextValue $ws.Range("D23") "8.490"
$ws.Range("E23").Value = "  -2.55%  "

Set-TextValue $ws.Range("D24") "5.078"
$ws.Range("E24").Value = "  -3.79%  "

Set-TextValue $ws.Range("D25") "138.04"
$ws.Range("E25").Value = "  -0.61%  "

Set-TextValue $ws.Range("D26") "15.27"
$ws.Range("E26").Value = "  -0.89%  "

Set-TextValue $ws.Range("D27") "1.438"
$ws.Range("E27").Value = "  +1.13%  "

Set-TextValue $ws.Range("D28") "1.744"
$ws.Range("E28").Value = "  -0.75%  "

Set-TextValue $ws.Range("D29") "105.92"
$ws.Range("E29").Value = "  -1.37%  "

Set-TextValue $ws.Range("D30") "3.915"
$ws.Range("E30").Value = "  -2.75%  "

Set-TextValue $ws.Range("D31") "0.07946"
$ws.Range("E31").Value = "  -1.21%  "

Set-TextValue $ws.Range("D32") "3.632"
$ws.Range("E32").Value = "  -2.59%  "

Set-TextValue $ws.Range("D33") "0.04513"
$ws.Range("E33").Value = "  -2.34%  "

Set-TextValue $ws.Range("D34") "2.626"
$ws.Range("E34").Value = "  +0.60%  "

Set-TextValue $ws.Range("D35") "0.9979"
$ws.Range("E35").Value = "  -1.61%  "

Set-TextValue $ws.Range("D36") "0.6232"
$ws.Range("E36").Value = "  -2.35%  "

Set-TextValue $ws.Range("D37") "0.9402"
$ws.Range("E37").Value = "  +3.98%  "

Set-TextValue $ws.Range("D38") "2.004"
$ws.Range("E38").Value = "  -2.98%  "

Set-TextValue $ws.Range("D39") "2.414"
$ws.Range("E39").Value = "  -0.52%  "

Set-TextValue $ws.Range("D40") "0.9982"
$ws.Range("E40").Value = "  -0.44%  "

Set-TextValue $ws.Range("D41") "0.01494"
$ws.Range("E41").Value = "  -0.71%  "

Set-TextValue $ws.Range("D42") "99.39"
$ws.Range("E42").Value = "  -2.52%  "

Set-TextValue $ws.Range("D43") "5.521"
$ws.Range("E43").Value = "  +1.63%  "

Set-TextValue $ws.Range("D44") "0.3822"
$ws.Range("E44").Value = "  -2.87%  "

Set-TextValue $ws.Range("D45") "6.908"
$ws.Range("E45").Value = "  +0.63%  "

Set-TextValue $ws.Range("D46") "0.1155"
$ws.Range("E46").Value = "  -2.54%  "

$ws.Range("E47").Value = "  +0.07%  "

Set-TextValue $ws.Range("D50") "51.46"
$ws.Range("E50").Value = "  -0.71%  "

Set-TextValue $ws.Range("D51") "1.218"
$ws.Range("E51").Value = "  -2.93%  "

# Rows 48/49: Elrond and EnergySwap swap order
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D48") "7.788"
$ws.Range("E48").Value = "  +0.18%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue $ws.Range("D49") "30.36"
$ws.Range("E49").Value = "  -0.82%  "

